$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '314.47'
    'E2' = '2.02%'
    'D3' = '39.23'
    'E3' = '-1.53%'
    'D4' = '5.146'
    'E4' = '-0.09%'
    'D5' = '0.08165'
    'E5' = '0.28%'
    'D6' = '1.993'
    'E6' = '2.67%'
    'D7' = '4.381'
    'E7' = '3.29%'
    'D8' = '8.343'
    'E8' = '2.16%'
    'D9' = '0.9365'
    'E9' = '0.76%'
    'D10' = '0.1304'
    'E10' = '-9.03%'
    'D11' = '0.1970'
    'E11' = '2.35%'
    'D12' = '0.08990'
    'E12' = '-1.07%'
    'D13' = '0.03529'
    'E13' = '0.38%'
    'D14' = '0.09741'
    'E14' = '-0.37%'
    'E15' = '0.88%'
    'D16' = '0.006590'
    'E16' = '12.98%'
    'D17' = '3.630'
    'E17' = '-7.41%'
    'D18' = '3.126'
    'E18' = '-7.52%'
    'D19' = '0.3470'
    'E19' = '1.15%'
    'D20' = '0.1309'
    'E20' = '-0.23%'
    'D21' = '5.001'
    'E21' = '7.47%'
    'D22' = '0.2490'
    'E22' = '2.72%'
    'D23' = '0.04362'
    'E23' = '-0.16%'
    'D24' = '0.001240'
    'E24' = '0.92%'
    'D25' = '0.004759'
    'E25' = '8.76%'
    'D26' = '0.0003895'
    'E26' = '199.55%'
    'E27' = '-7.57%'
    'D39' = '0.02232'
    'E39' = '8.66%'
    'D40' = '0.05186'
    'E40' = '2.27%'
    'D41' = '0.007760'
    'E41' = '5.14%'
    'D42' = '0.01031'
    'E42' = '4.58%'
    'D43' = '0.1398'
    'E43' = '2.53%'
    'D44' = '0.002103'
    'E44' = '-1.30%'
    'D45' = '0.008859'
    'E45' = '-5.54%'
    'D46' = '0.00006823'
    'E46' = '7.13%'
    'D47' = '0.00000000751'
    'E47' = '0.07%'
    'D48' = '0.003009'
    'E48' = '10.88%'
    'D49' = '0.001692'
    'E49' = '30.07%'
    'D50' = '0.00002102'
    'E50' = '0.07%'
    'D51' = '0.0002002'
    'E51' = '0.07%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
